# Update crypto price (D) and volume-change (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '25.573.18'
$ws.Cells.Item(2, 5).Value = '  +2.60%  '
$ws.Cells.Item(3, 4).Value = '1.669.90'
$ws.Cells.Item(3, 5).Value = '  +2.07%  '
$ws.Cells.Item(4, 4).Value = '''0.9992'
$ws.Cells.Item(4, 5).Value = '  +0.08%  '
$ws.Cells.Item(5, 4).Value = '''239.25'
$ws.Cells.Item(5, 5).Value = '  +1.63%  '
$ws.Cells.Item(6, 4).Value = '''1.000'
$ws.Cells.Item(6, 5).Value = '  -0.09%  '
$ws.Cells.Item(7, 4).Value = '''0.4796'
$ws.Cells.Item(7, 5).Value = '  +1.56%  '
$ws.Cells.Item(8, 4).Value = '''0.2629'
$ws.Cells.Item(8, 5).Value = '  +2.94%  '
$ws.Cells.Item(9, 4).Value = '''0.06170'
$ws.Cells.Item(9, 5).Value = '  +2.93%  '
$ws.Cells.Item(10, 4).Value = '1.668.15'
$ws.Cells.Item(10, 5).Value = '  +2.00%  '
$ws.Cells.Item(11, 4).Value = '''0.06990'
$ws.Cells.Item(11, 5).Value = '  -2.48%  '
$ws.Cells.Item(12, 5).Value = '  +1.00%  '
$ws.Cells.Item(13, 4).Value = '''0.5885'
$ws.Cells.Item(13, 5).Value = '  -4.00%  '
$ws.Cells.Item(14, 4).Value = '''4.382'
$ws.Cells.Item(14, 5).Value = '  -0.58%  '
$ws.Cells.Item(15, 5).Value = '  +3.99%  '
$ws.Cells.Item(16, 4).Value = '''1.0000'
$ws.Cells.Item(16, 5).Value = '  -0.22%  '
$ws.Cells.Item(17, 4).Value = '''0.9999'
$ws.Cells.Item(17, 5).Value = '  +0.17%  '
$ws.Cells.Item(18, 4).Value = '25.561.23'
$ws.Cells.Item(18, 5).Value = '  +2.64%  '
$ws.Cells.Item(19, 4).Value = '''0.000006758'
$ws.Cells.Item(19, 5).Value = '  +3.03%  '
$ws.Cells.Item(20, 4).Value = '''11.45'
$ws.Cells.Item(20, 5).Value = '  +2.34%  '
$ws.Cells.Item(21, 4).Value = '1.884.16'
$ws.Cells.Item(21, 5).Value = '  +1.93%  '
$ws.Cells.Item(22, 4).Value = '''4.436'
$ws.Cells.Item(22, 5).Value = '  +0.82%  '
$ws.Cells.Item(23, 5).Value = '  +2.21%  '
$ws.Cells.Item(24, 4).Value = '''5.276'
$ws.Cells.Item(24, 5).Value = '  +0.56%  '
$ws.Cells.Item(25, 4).Value = '''136.71'
$ws.Cells.Item(25, 5).Value = '  +3.49%  '
$ws.Cells.Item(26, 4).Value = '''15.05'
$ws.Cells.Item(26, 5).Value = '  +1.92%  '
$ws.Cells.Item(27, 5).Value = '  +1.34%  '
$ws.Cells.Item(28, 4).Value = '''1.720'
$ws.Cells.Item(28, 5).Value = '  +4.34%  '
$ws.Cells.Item(29, 4).Value = '''104.82'
$ws.Cells.Item(29, 5).Value = '  +1.99%  '
$ws.Cells.Item(30, 4).Value = '''3.968'
$ws.Cells.Item(30, 5).Value = '  +7.03%  '
$ws.Cells.Item(31, 4).Value = '''0.07830'
$ws.Cells.Item(31, 5).Value = '  +1.11%  '
$ws.Cells.Item(32, 4).Value = '''3.648'
$ws.Cells.Item(32, 5).Value = '  +3.27%  '
$ws.Cells.Item(33, 5).Value = '  -0.15%  '
$ws.Cells.Item(34, 4).Value = '''0.04225'
$ws.Cells.Item(34, 5).Value = '  -3.35%  '
$ws.Cells.Item(35, 4).Value = '''2.619'
$ws.Cells.Item(35, 5).Value = '  +0.80%  '
$ws.Cells.Item(36, 4).Value = '''0.6094'
$ws.Cells.Item(36, 5).Value = '  +5.15%  '
$ws.Cells.Item(37, 4).Value = '''0.9535'
$ws.Cells.Item(37, 5).Value = '  +4.08%  '
$ws.Cells.Item(38, 5).Value = '  +3.31%  '
$ws.Cells.Item(39, 4).Value = '''0.8599'
$ws.Cells.Item(39, 5).Value = '  +4.43%  '
$ws.Cells.Item(40, 4).Value = '''0.9994'
$ws.Cells.Item(40, 5).Value = '  +0.08%  '
$ws.Cells.Item(41, 5).Value = '  +4.85%  '
$ws.Cells.Item(42, 4).Value = '''0.01476'
$ws.Cells.Item(42, 5).Value = '  -4.99%  '
$ws.Cells.Item(43, 4).Value = '''96.53'
$ws.Cells.Item(43, 5).Value = '  -0.72%  '
$ws.Cells.Item(44, 4).Value = '''0.3766'
$ws.Cells.Item(44, 5).Value = '  +1.95%  '
$ws.Cells.Item(45, 4).Value = '''4.869'
$ws.Cells.Item(45, 5).Value = '  +2.85%  '
$ws.Cells.Item(46, 4).Value = '''0.1118'
$ws.Cells.Item(46, 5).Value = '  -1.64%  '
$ws.Cells.Item(47, 4).Value = '''6.221'
$ws.Cells.Item(47, 5).Value = '  +2.51%  '
$ws.Cells.Item(48, 4).Value = '''0.05261'
$ws.Cells.Item(48, 5).Value = '  +1.34%  '
$ws.Cells.Item(49, 5).Value = '  +1.41%  '
$ws.Cells.Item(50, 4).Value = '''7.381'
$ws.Cells.Item(50, 5).Value = '  +3.15%  '
$ws.Cells.Item(51, 4).Value = '''1.001'
$ws.Cells.Item(51, 5).Value = '  +0.16%  '
